$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style (bold, centered, bordered) from existing header cell E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean data rows 2-12 for columns F (KNN), G (SVM), H (RF)
$values = @{
    2  = @($true, $false, $true)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $false, $false)
    6  = @($false, $false, $false)
    7  = @($true, $false, $false)
    8  = @($true, $true, $true)
    9  = @($false, $false, $false)
    10 = @($false, $false, $false)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
}

foreach ($row in $values.Keys) {
    $trio = $values[$row]
    $ws.Cells.Item($row, 6).Value = $trio[0]
    $ws.Cells.Item($row, 7).Value = $trio[1]
    $ws.Cells.Item($row, 8).Value = $trio[2]
}
